# Apply "Exp 12 Parameters" update to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing typo in F7: "Exp11.png" -> "Exp 11.png"
$ws.Range("F7").Value = "Exp 11.png"

# Fill in the new row 8 with the Exp 12 parameters
$ws.Range("A8").Value = "Exp 12"
$ws.Range("B8").Value = 0.8
$ws.Range("C8").Value = 1
$ws.Range("F8").Value = "Exp 12.png"

# Move the active selection to F11 (next empty row to fill)
$ws.Range("F11").Select()

$wb.Save()
